$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# Row 27: add a new status value "?" in column E
$ws.Range("E27").Value = "?"

# Row 28: update run timestamp in column B, clear the error status in column E
$ws.Range("B28").Value = "2022-06-15 16-04-57"
$ws.Range("E28").ClearContents()

# Row 30: update run timestamp in column B, clear the error status in column E
$ws.Range("B30").Value = "2022-06-15 16-05-31"
$ws.Range("E30").ClearContents()

# Move the active selection to E30
$ws.Range("E30").Select()
